$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.379.92'
$ws.Range("E2").Value = '  +1.89%  '
$ws.Range("D3").Value = '2.595.14'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '569.60'
$ws.Range("E5").Value = '  +1.73%  '
$ws.Range("D6").Value = '141.96'
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '0.600'
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("D9").Value = '2.615.66'
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("D10").Value = '6.54'
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("E12").Value = '  +2.69%  '
$ws.Range("D13").Value = '0.149'
$ws.Range("E13").Value = '  -6.05%  '
$ws.Range("D14").Value = '3.061.36'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '60.385.46'
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").Value = '23.26'
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("E17").Value = '  +2.75%  '
$ws.Range("D18").Value = '2.607.92'
$ws.Range("E18").Value = '  +1.24%  '
$ws.Range("E19").Value = '  +9.19%  '
$ws.Range("E20").Value = '  +1.38%  '
$ws.Range("D21").Value = '345.99'
$ws.Range("E21").Value = '  +2.70%  '
$ws.Range("D22").Value = '6.99'
$ws.Range("E22").Value = '  +8.92%  '
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("D24").Value = '0.527'
$ws.Range("E24").Value = '  +12.88%  '
$ws.Range("D25").Value = '63.33'
$ws.Range("E25").Value = '  -1.09%  '
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("E27").Value = '  -1.86%  '
$ws.Range("D28").Value = '7.69'
$ws.Range("E28").Value = '  +4.46%  '
$ws.Range("D29").Value = '0.0₃0782'
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("D30").Value = '1.81'
$ws.Range("E30").Value = '  +8.84%  '
$ws.Range("E31").Value = '  +4.00%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").Value = '160.30'
$ws.Range("E33").Value = '  +0.75%  '
$ws.Range("E34").Value = '  +2.38%  '
$ws.Range("E35").Value = '  +4.60%  '
$ws.Range("D36").Value = '0.967'
$ws.Range("E36").Value = '  +9.75%  '
$ws.Range("E37").Value = '  +3.75%  '
$ws.Range("E38").Value = '  +8.49%  '
$ws.Range("D39").Value = '37.66'
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("E40").Value = '  +3.43%  '
$ws.Range("D41").Value = '0.848'
$ws.Range("E41").Value = '  -2.22%  '
$ws.Range("D42").Value = '294.00'
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("D43").Value = '138.82'
$ws.Range("E43").Value = '  +4.68%  '
$ws.Range("E44").Value = '  -0.29%  '
$ws.Range("D45").Value = '0.0983'
$ws.Range("E45").Value = '  +0.79%  '
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("D47").Value = '19.66'
$ws.Range("E47").Value = '  +3.30%  '
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("E49").Value = '  +2.35%  '
$ws.Range("D50").Value = '19.83'
$ws.Range("E50").Value = '  +6.33%  '
$ws.Range("D51").Value = '10.72'
$ws.Range("E51").Value = '  +0.81%  '
